$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Week 1")
$ws.Range("A1").Value = "TEST"
